# Auto-generated Excel COM-interop edit script
# Línea 141 horarios update: 07:21:42 -> 07:48:14 scrape, with resorted + appended rows
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 07:48:14"
$ws1.Range("A3").Value = "Total filas: 80"
$ws1.Cells.Item(18,1).Value = "05:23:04"
$ws1.Cells.Item(18,3).Value = "10_OLMOS"
$ws1.Cells.Item(18,4).Value = 41
$ws1.Cells.Item(19,1).Value = "04:56:49"
$ws1.Cells.Item(19,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(19,4).Value = 68
$ws1.Cells.Item(28,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(29,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(55,1).Value = "07:48:14"
$ws1.Cells.Item(55,2).Value = "07:48"
$ws1.Cells.Item(55,3).Value = "215A_EL PATO"
$ws1.Cells.Item(55,4).Value = 0
$ws1.Cells.Item(56,2).Value = "07:49"
$ws1.Cells.Item(56,3).Value = "15_ABASTO"
$ws1.Cells.Item(56,4).Value = 28
$ws1.Cells.Item(57,1).Value = "07:48:14"
$ws1.Cells.Item(57,2).Value = "07:58"
$ws1.Cells.Item(57,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(57,4).Value = 10
$ws1.Cells.Item(58,1).Value = "07:21:42"
$ws1.Cells.Item(58,2).Value = "07:59"
$ws1.Cells.Item(58,4).Value = 38
$ws1.Cells.Item(59,1).Value = "07:48:14"
$ws1.Cells.Item(59,2).Value = "08:00"
$ws1.Cells.Item(59,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(59,4).Value = 12
$ws1.Cells.Item(60,1).Value = "07:48:14"
$ws1.Cells.Item(60,2).Value = "08:01"
$ws1.Cells.Item(60,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(60,4).Value = 13
$ws1.Cells.Item(61,1).Value = "07:48:14"
$ws1.Cells.Item(61,2).Value = "08:03"
$ws1.Cells.Item(61,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(61,4).Value = 15
$ws1.Cells.Item(62,1).Value = "06:46:06"
$ws1.Cells.Item(62,2).Value = "08:03"
$ws1.Cells.Item(62,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(62,4).Value = 77
$ws1.Cells.Item(63,1).Value = "06:58:01"
$ws1.Cells.Item(63,2).Value = "08:04"
$ws1.Cells.Item(63,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(63,4).Value = 66
$ws1.Cells.Item(64,1).Value = "07:48:14"
$ws1.Cells.Item(64,2).Value = "08:14"
$ws1.Cells.Item(64,3).Value = "10_OLMOS"
$ws1.Cells.Item(64,4).Value = 26
$ws1.Cells.Item(65,1).Value = "07:48:14"
$ws1.Cells.Item(65,2).Value = "08:19"
$ws1.Cells.Item(65,3).Value = "15_ABASTO"
$ws1.Cells.Item(65,4).Value = 31
$ws1.Cells.Item(66,1).Value = "07:48:14"
$ws1.Cells.Item(66,2).Value = "08:21"
$ws1.Cells.Item(66,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(66,4).Value = 33
$ws1.Cells.Item(67,1).Value = "07:21:42"
$ws1.Cells.Item(67,2).Value = "08:29"
$ws1.Cells.Item(67,3).Value = "14_ABASTO"
$ws1.Cells.Item(67,4).Value = 68
$ws1.Cells.Item(68,1).Value = "07:48:14"
$ws1.Cells.Item(68,2).Value = "08:30"
$ws1.Cells.Item(68,3).Value = "14_ABASTO"
$ws1.Cells.Item(68,4).Value = 42
$ws1.Cells.Item(69,2).Value = "08:33"
$ws1.Cells.Item(69,3).Value = "215C_EL PATO"
$ws1.Cells.Item(69,4).Value = 72
$ws1.Cells.Item(70,1).Value = "07:48:14"
$ws1.Cells.Item(70,2).Value = "08:34"
$ws1.Cells.Item(70,3).Value = "215C_EL PATO"
$ws1.Cells.Item(70,4).Value = 46
$ws1.Cells.Item(71,1).Value = "07:48:14"
$ws1.Cells.Item(71,2).Value = "08:48"
$ws1.Cells.Item(71,3).Value = "215A_EL PATO"
$ws1.Cells.Item(71,4).Value = 60
$ws1.Cells.Item(72,2).Value = "08:51"
$ws1.Cells.Item(72,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(72,4).Value = 90
$ws1.Cells.Item(73,1).Value = "07:48:14"
$ws1.Cells.Item(73,2).Value = "08:52"
$ws1.Cells.Item(73,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(73,4).Value = 64
$ws1.Cells.Item(74,1).Value = "07:21:42"
$ws1.Cells.Item(74,2).Value = "08:59"
$ws1.Cells.Item(74,3).Value = "215B_EL PATO"
$ws1.Cells.Item(74,4).Value = 98
$ws1.Cells.Item(74,5).Value = "LP1912"
$ws1.Cells.Item(75,1).Value = "07:48:14"
$ws1.Cells.Item(75,2).Value = "09:00"
$ws1.Cells.Item(75,3).Value = "215B_EL PATO"
$ws1.Cells.Item(75,4).Value = 72
$ws1.Cells.Item(75,5).Value = "LP1912"
$ws1.Cells.Item(76,1).Value = "07:48:14"
$ws1.Cells.Item(76,2).Value = "09:03"
$ws1.Cells.Item(76,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(76,4).Value = 75
$ws1.Cells.Item(76,5).Value = "LP1912"
$ws1.Cells.Item(77,1).Value = "07:48:14"
$ws1.Cells.Item(77,2).Value = "09:03"
$ws1.Cells.Item(77,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(77,4).Value = 75
$ws1.Cells.Item(77,5).Value = "LP1912"
$ws1.Cells.Item(78,1).Value = "07:21:42"
$ws1.Cells.Item(78,2).Value = "09:14"
$ws1.Cells.Item(78,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(78,4).Value = 113
$ws1.Cells.Item(78,5).Value = "LP1912"
$ws1.Cells.Item(79,1).Value = "07:48:14"
$ws1.Cells.Item(79,2).Value = "09:15"
$ws1.Cells.Item(79,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(79,4).Value = 87
$ws1.Cells.Item(79,5).Value = "LP1912"
$ws1.Cells.Item(80,1).Value = "07:48:14"
$ws1.Cells.Item(80,2).Value = "09:16"
$ws1.Cells.Item(80,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(80,4).Value = 88
$ws1.Cells.Item(80,5).Value = "LP1912"
$ws1.Cells.Item(81,1).Value = "07:21:42"
$ws1.Cells.Item(81,2).Value = "09:18"
$ws1.Cells.Item(81,3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(81,4).Value = 117
$ws1.Cells.Item(81,5).Value = "LP1912"
$ws1.Cells.Item(82,1).Value = "07:48:14"
$ws1.Cells.Item(82,2).Value = "09:19"
$ws1.Cells.Item(82,3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(82,4).Value = 91
$ws1.Cells.Item(82,5).Value = "LP1912"
$ws1.Cells.Item(83,1).Value = "07:48:14"
$ws1.Cells.Item(83,2).Value = "09:29"
$ws1.Cells.Item(83,3).Value = "10_OLMOS"
$ws1.Cells.Item(83,4).Value = 101
$ws1.Cells.Item(83,5).Value = "LP1912"
$ws1.Cells.Item(84,1).Value = "07:48:14"
$ws1.Cells.Item(84,2).Value = "09:34"
$ws1.Cells.Item(84,3).Value = "15_ABASTO"
$ws1.Cells.Item(84,4).Value = 106
$ws1.Cells.Item(84,5).Value = "LP1912"
$ws1.Cells.Item(85,1).Value = "07:48:14"
$ws1.Cells.Item(85,2).Value = "09:45"
$ws1.Cells.Item(85,3).Value = "14_ABASTO"
$ws1.Cells.Item(85,4).Value = 117
$ws1.Cells.Item(85,5).Value = "LP1912"

$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 07:48:14"
$ws2.Range("A3").Value = "Total filas: 18"
$ws2.Cells.Item(16,1).Value = "07:48:14"
$ws2.Cells.Item(16,2).Value = "07:48"
$ws2.Cells.Item(16,3).Value = "215A_EL PATO"
$ws2.Cells.Item(16,4).Value = 0
$ws2.Cells.Item(17,1).Value = "07:21:42"
$ws2.Cells.Item(17,2).Value = "08:33"
$ws2.Cells.Item(17,4).Value = 72
$ws2.Cells.Item(18,1).Value = "07:48:14"
$ws2.Cells.Item(18,2).Value = "08:34"
$ws2.Cells.Item(18,3).Value = "215C_EL PATO"
$ws2.Cells.Item(18,4).Value = 46
$ws2.Cells.Item(19,1).Value = "07:48:14"
$ws2.Cells.Item(19,2).Value = "08:48"
$ws2.Cells.Item(19,3).Value = "215A_EL PATO"
$ws2.Cells.Item(19,4).Value = 60
$ws2.Cells.Item(20,2).Value = "08:59"
$ws2.Cells.Item(20,3).Value = "215B_EL PATO"
$ws2.Cells.Item(20,4).Value = 98
$ws2.Cells.Item(21,1).Value = "07:48:14"
$ws2.Cells.Item(21,2).Value = "09:00"
$ws2.Cells.Item(21,3).Value = "215B_EL PATO"
$ws2.Cells.Item(21,4).Value = 72
$ws2.Cells.Item(21,5).Value = "LP1912"
$ws2.Cells.Item(22,1).Value = "07:21:42"
$ws2.Cells.Item(22,2).Value = "09:18"
$ws2.Cells.Item(22,3).Value = "215_EL PELIGRO"
$ws2.Cells.Item(22,4).Value = 117
$ws2.Cells.Item(22,5).Value = "LP1912"
$ws2.Cells.Item(23,1).Value = "07:48:14"
$ws2.Cells.Item(23,2).Value = "09:19"
$ws2.Cells.Item(23,3).Value = "215_EL PELIGRO"
$ws2.Cells.Item(23,4).Value = 91
$ws2.Cells.Item(23,5).Value = "LP1912"

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 07:48:14"
$ws3.Range("A3").Value = "Total filas: 13"
$ws3.Cells.Item(14,1).Value = "07:48:14"
$ws3.Cells.Item(14,4).Value = 22
$ws3.Cells.Item(16,1).Value = "07:48:14"
$ws3.Cells.Item(16,2).Value = "08:26"
$ws3.Cells.Item(16,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(16,4).Value = 38
$ws3.Cells.Item(16,5).Value = "L6203"
$ws3.Cells.Item(17,1).Value = "06:58:01"
$ws3.Cells.Item(17,2).Value = "08:52"
$ws3.Cells.Item(17,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(17,4).Value = 114
$ws3.Cells.Item(17,5).Value = "L6173"
$ws3.Cells.Item(18,1).Value = "07:48:14"
$ws3.Cells.Item(18,2).Value = "09:09"
$ws3.Cells.Item(18,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(18,4).Value = 81
$ws3.Cells.Item(18,5).Value = "L6173"

